$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: date 2021-08-21, 8 hours, activity text
$ws.Range("A26").Value = 44429
$ws.Range("B26").Value = 8
$ws.Range("D26").Value = "Implementierung des sendens von Beispieldaten des MeasuredData-Signals im RepitayaStub"

# Row 27: date 2021-08-22, 8 hours, activity text
$ws.Range("A27").Value = 44430
$ws.Range("B27").Value = 8
$ws.Range("D27").Value = "Implementierung des Empfangens der MeasuredData in der UI und Anzeige via der uPlot-Chart"

# Update selection to D28 to match the final cursor position
$ws.Range("D28").Select()
